$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4180726.5
$ws.Range("J17").Value = 4282691.5
$ws.Range("L17").Value = 12848074.5
$ws.Range("N17").Value = -12848410.5

$ws.Range("H70").Value = 2625.318
$ws.Range("J70").Value = 2483.75
$ws.Range("L70").Value = 7451.25
$ws.Range("N70").Value = -7991.25

$ws.Range("H73").Value = 2625.318
$ws.Range("J73").Value = 2483.75
$ws.Range("L73").Value = 7451.25
$ws.Range("N73").Value = -9323.25

$ws.Range("H129").Value = 864.6799999999999
$ws.Range("I129").Value = 464.85715
$ws.Range("J129").Value = 970.96204
$ws.Range("K129").Value = 1394.57145
$ws.Range("L129").Value = 2912.88612
$ws.Range("M129").Value = 3605.42855
$ws.Range("N129").Value = -12912.88612

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1193.3478
$ws.Range("I2").Value = 723.5263
$ws.Range("J2").Value = 3425
$ws.Range("K2").Value = 723.5263
$ws.Range("L2").Value = 3425
$ws.Range("M2").Value = -610.5263
$ws.Range("N2").Value = -3651

$ws.Range("H32").Value = 5359.4326
$ws.Range("I32").Value = 3267.3389
$ws.Range("J32").Value = 13588.333
$ws.Range("K32").Value = 3267.3389
$ws.Range("L32").Value = 13588.333
$ws.Range("M32").Value = -2980.3389
$ws.Range("N32").Value = -14162.333

$ws.Range("H116").Value = 1193.3478
$ws.Range("I116").Value = 723.5263
$ws.Range("J116").Value = 3425
$ws.Range("K116").Value = 723.5263
$ws.Range("L116").Value = 3425
$ws.Range("M116").Value = 1570.4737
$ws.Range("N116").Value = -8013

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1193.3478
$ws.Range("I3").Value = 723.5263
$ws.Range("J3").Value = 3425
$ws.Range("K3").Value = 723.5263
$ws.Range("L3").Value = 3425
$ws.Range("M3").Value = -609.5263
$ws.Range("N3").Value = -3653

$ws.Range("H134").Value = 1484.4133
$ws.Range("I134").Value = 1218.6936
$ws.Range("K134").Value = 3656.0808
$ws.Range("M134").Value = -1121.0808

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2404957
$ws.Range("I16").Value = 4274317
$ws.Range("J16").Value = 1494.5
$ws.Range("K16").Value = 4274317
$ws.Range("L16").Value = 1494.5
$ws.Range("M16").Value = -4274030
$ws.Range("N16").Value = -2068.5

$ws.Range("H31").Value = 3182.9487
$ws.Range("I31").Value = 1757.3462
$ws.Range("J31").Value = 6034.154
$ws.Range("K31").Value = 1757.3462
$ws.Range("L31").Value = 6034.154
$ws.Range("M31").Value = -1462.3462
$ws.Range("N31").Value = -6624.154

$ws.Range("H34").Value = 3182.9487
$ws.Range("I34").Value = 1757.3462
$ws.Range("J34").Value = 6034.154
$ws.Range("K34").Value = 1757.3462
$ws.Range("L34").Value = 6034.154
$ws.Range("M34").Value = -1555.3462
$ws.Range("N34").Value = -6438.154

$ws.Range("H55").Value = 15000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H113").Value = 2404957
$ws.Range("I113").Value = 4274317
$ws.Range("J113").Value = 1494.5
$ws.Range("K113").Value = 4274317
$ws.Range("L113").Value = 1494.5
$ws.Range("M113").Value = -4272147
$ws.Range("N113").Value = -5834.5

$ws.Range("H132").Value = 1687.4694
$ws.Range("I132").Value = 1466.1428
$ws.Range("J132").Value = 3015.4285
$ws.Range("K132").Value = 4398.428400000001
$ws.Range("L132").Value = 9046.2855
$ws.Range("M132").Value = -1868.428400000001
$ws.Range("N132").Value = -14106.2855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 262022.3
$ws.Range("I5").Value = 773.4545000000001
$ws.Range("J5").Value = 501500.4
$ws.Range("K5").Value = 2320.3635
$ws.Range("L5").Value = 1504501.2
$ws.Range("M5").Value = -2208.3635
$ws.Range("N5").Value = -1504725.2

$ws.Range("H122").Value = 583.4737
$ws.Range("I122").Value = 388.36365
$ws.Range("J122").Value = 851.75
$ws.Range("K122").Value = 3495.27285
$ws.Range("L122").Value = 7665.75
$ws.Range("M122").Value = -1045.27285
$ws.Range("N122").Value = -12565.75

$ws.Range("H131").Value = 1640327.5
$ws.Range("J131").Value = 1060.1296
$ws.Range("L131").Value = 3180.3888
$ws.Range("N131").Value = -13260.3888

$ws.Range("H132").Value = 1331.8182
$ws.Range("I132").Value = 864.2857
$ws.Range("J132").Value = 2150
$ws.Range("K132").Value = 7778.571300000001
$ws.Range("L132").Value = 19350
$ws.Range("M132").Value = -5248.571300000001
$ws.Range("N132").Value = -24410

$ws.Range("H135").Value = 262022.3
$ws.Range("I135").Value = 773.4545000000001
$ws.Range("J135").Value = 501500.4
$ws.Range("K135").Value = 6961.0905
$ws.Range("L135").Value = 4513503.600000001
$ws.Range("M135").Value = -4426.0905
$ws.Range("N135").Value = -4518573.600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1478.6154
$ws.Range("I102").Value = 1215.25
$ws.Range("K102").Value = 1215.25
$ws.Range("M102").Value = 406.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2220
$ws.Range("I7").Value = 2128.5715
$ws.Range("J7").Value = 3500
$ws.Range("K7").Value = 2128.5715
$ws.Range("L7").Value = 3500
$ws.Range("M7").Value = -2016.5715
$ws.Range("N7").Value = -3724

$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H126").Value = 2220
$ws.Range("I126").Value = 2128.5715
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 6385.7145
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -3915.7145
$ws.Range("N126").Value = -15440

$ws.Range("H136").Value = 3718.6611
$ws.Range("I136").Value = 1632.3877
$ws.Range("J136").Value = 13941.4
$ws.Range("K136").Value = 4897.1631
$ws.Range("L136").Value = 41824.2
$ws.Range("M136").Value = -2347.1631
$ws.Range("N136").Value = -46924.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4071.4285
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 4071.4285
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 8142.857
$ws.Range("N81").Value = -10264.857
$ws.Range("M81").ClearContents()

$ws.Range("H84").Value = 4071.4285
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 4071.4285
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 40714.285
$ws.Range("N84").Value = -51322.285
$ws.Range("M84").ClearContents()

$ws.Range("H132").Value = 1138.8235
$ws.Range("I132").Value = 834.425
$ws.Range("J132").Value = 2245.7273
$ws.Range("K132").Value = 2503.275
$ws.Range("L132").Value = 6737.1819
$ws.Range("M132").Value = 26.72500000000036
$ws.Range("N132").Value = -11797.1819

$ws.Range("H136").Value = 840.3269
$ws.Range("I136").Value = 432.2973
$ws.Range("J136").Value = 1846.8
$ws.Range("K136").Value = 1296.8919
$ws.Range("L136").Value = 5540.4
$ws.Range("M136").Value = 1253.1081
$ws.Range("N136").Value = -10640.4
